# Auto-generated script to apply cryptos.xlsx updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.139.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.356.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.643"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.896.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.398.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.959.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +5.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "303.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.55%  "
$ws.Range("E29").Value = "  +6.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.77%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0507"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.68%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.183.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.42%  "
